$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: new "common part" summary row (m0004)
$ws.Range("A4").Value = 'm0004'
$ws.Range("B4").Value = '지금까지의 각 참/거짓 결과에서 ㄱ, ㄴ, ㄷ 중 옳은 것을 모두 고릅니다.'

# Make room for the new d0018-d0021 block: shift rows 61/71/81 down to 71/81/91
$ws.Rows("61:70").Insert()

# Row 58: d0018
$ws.Range("A58").Value = 'd0018'
$ws.Range("B58").Value = '알아 낼 수 있는 위치, 즉 $x$ 좌표를 각각의 시각에 대해 알아냅니다.'
$ws.Range("C58").Value = '$x(0)$과 $x(1)$; '

# Row 59: d0019
$ws.Range("A59").Value = 'd0019'
$ws.Range("B59").Value = '두 시각 사이의 위치의 변화량, 즉 $x$ 좌표의 변화량을 정적분으로 알아냅니다.'
$ws.Range("C59").Value = '$x(1) - x(0)$;'

# Row 60: d0020
$ws.Range("A60").Value = 'd0020'
$ws.Range("B60").Value = '범위내의 모든 시각에 대한 위치의 주어진 조건과 문제 전체의 전제 조건인 실제 움직인 거리와 비교해 봅니다.'
$ws.Range("C60").Value = '$\left|x\left(t_{1}\right)\right|>1$;'

# Row 61: d0021
$ws.Range("A61").Value = 'd0021'
$ws.Range("B61").Value = '두 시각 사이의 어떤 위치의 조건과 문제 전체의 전제 조건인 실제 움직인 거리와 비교해 봅니다.'
$ws.Range("C61").Value = '$|x(t)| < 1$;'

[void]$ws.Range("B4").Select()
